{"js": "// Replace the essay's title, author byline, author e-mail, body copy and\n// summary with the new \"Chemistry\" themed content, and add a trailing\n// empty paragraph at the end of the document.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 0: Title\nitems[0].insertText(\n  \"The Marvelous World of Chemistry: Unraveling the Secrets of Matter\",\n  Word.InsertLocation.replace\n);\n\n// 1: Author name (\"Richard T. Pomeroy\" -> \"Dr. Evelyn Richards\")\nitems[1].insertText(\"Dr. Evelyn Richards\", Word.InsertLocation.replace);\n\n// 2: Author e-mail (\"richard.pomeroy@biotechne.com\" -> \"evrichards@eduworld.org\")\nitems[2].insertText(\"evrichards@eduworld.org\", Word.InsertLocation.replace);\n\n// 3: blank separator paragraph stays blank - no change needed.\n\n// 4: Main body paragraph (three \"paragraphs\" separated by manual line breaks)\nconst bodyText =\n  \"In the vast expanse of human knowledge, chemistry stands tall as a beacon of understanding, illuminating the innermost workings of matter and its myriad transformations.\" +\n  \" It is a voyage into the heart of creation, a quest to unlock the secrets that govern the composition, structure, and behavior of substances that make up our world.\" +\n  \" As we delve into the remarkable field of chemistry, we embark on a journey of discovery, where each step brings us closer to comprehending the intricate symphony of interactions that orchestrate the universe around us.\" +\n  \"\\u000b\\u000bFrom the simplest elements to the most complex compounds, chemistry unveils the fundamental building blocks of our existence.\" +\n  \" Within the atoms and molecules, we find a hidden world of particles, each possessing unique characteristics and intricate relationships.\" +\n  \" Through careful experimentation and theoretical exploration, chemists strive to unravel the intricate tapestry of these interactions, seeking to understand not only the nature of matter itself but also the forces that shape its behavior.\" +\n  \"\\u000b\\u000bChemistry plays a pivotal role in defining the world we inhabit, influencing everything from the air we breathe to the food we consume.\" +\n  \" It shapes the properties of materials, determining their strength, flexibility, and reactivity.\" +\n  \" It governs the interactions between living organisms, dictating the intricate web of life's processes.\" +\n  \" Chemistry holds the key to understanding the delicate balance of our planet, providing insights into the challenges of pollution, climate change, and the sustainable use of resources.\";\nitems[4].insertText(bodyText, Word.InsertLocation.replace);\n\n// 5: \"Summary\" heading stays the same - no change needed.\n\n// 6: Summary body paragraph\nconst summaryText =\n  \"Chemistry is the science that explores the composition, structure, and behavior of matter.\" +\n  \" It delves into the innermost workings of substances, seeking to understand the nature of atoms, molecules, and their interactions.\" +\n  \" Chemistry plays a fundamental role in shaping the world we inhabit, influencing everything from materials science to life processes.\" +\n  \" It empowers us to understand the delicate balance of our planet and address critical challenges such as pollution and climate change.\" +\n  \" By unraveling the secrets of matter, chemistry becomes an invaluable tool in advancing human knowledge and shaping a sustainable future.\";\nitems[6].insertText(summaryText, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// Append a new empty paragraph at the very end of the document.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Replace the essay's title, author byline, author e-mail, body copy and\n# summary with the new \"Chemistry\" themed content, and add a trailing\n# empty paragraph at the end of the document.\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($paragraph, [string]$newText) {\n    # Paragraph.Range includes the trailing paragraph mark; replacing that\n    # whole range (instead of just assigning Range.Text) collapses the\n    # paragraph down to a single run carrying the first run's formatting\n    # and swaps in the new text, same as Office.js insertText(...,replace).\n    $r = $paragraph.Range\n    $textRange = $d.Range($r.Start, $r.End - 1)\n    $textRange.Text = $newText\n}\n\n# 1: Title\nSet-ParagraphText $d.Paragraphs.Item(1) \"The Marvelous World of Chemistry: Unraveling the Secrets of Matter\"\n\n# 2: Author name (\"Richard T. Pomeroy\" -> \"Dr. Evelyn Richards\")\nSet-ParagraphText $d.Paragraphs.Item(2) \"Dr. Evelyn Richards\"\n\n# 3: Author e-mail (\"richard.pomeroy@biotechne.com\" -> \"evrichards@eduworld.org\")\nSet-ParagraphText $d.Paragraphs.Item(3) \"evrichards@eduworld.org\"\n\n# 4: blank separator paragraph stays blank - no change needed.\n\n# 5: Main body paragraph (three \"paragraphs\" separated by manual line breaks)\n$bodyText = \"In the vast expanse of human knowledge, chemistry stands tall as a beacon of understanding, illuminating the innermost workings of matter and its myriad transformations.\" + `\n  \" It is a voyage into the heart of creation, a quest to unlock the secrets that govern the composition, structure, and behavior of substances that make up our world.\" + `\n  \" As we delve into the remarkable field of chemistry, we embark on a journey of discovery, where each step brings us closer to comprehending the intricate symphony of interactions that orchestrate the universe around us.\" + `\n  \"`v`vFrom the simplest elements to the most complex compounds, chemistry unveils the fundamental building blocks of our existence.\" + `\n  \" Within the atoms and molecules, we find a hidden world of particles, each possessing unique characteristics and intricate relationships.\" + `\n  \" Through careful experimentation and theoretical exploration, chemists strive to unravel the intricate tapestry of these interactions, seeking to understand not only the nature of matter itself but also the forces that shape its behavior.\" + `\n  \"`v`vChemistry plays a pivotal role in defining the world we inhabit, influencing everything from the air we breathe to the food we consume.\" + `\n  \" It shapes the properties of materials, determining their strength, flexibility, and reactivity.\" + `\n  \" It governs the interactions between living organisms, dictating the intricate web of life's processes.\" + `\n  \" Chemistry holds the key to understanding the delicate balance of our planet, providing insights into the challenges of pollution, climate change, and the sustainable use of resources.\"\nSet-ParagraphText $d.Paragraphs.Item(5) $bodyText\n\n# 6: \"Summary\" heading stays the same - no change needed.\n\n# 7: Summary body paragraph\n$summaryText = \"Chemistry is the science that explores the composition, structure, and behavior of matter.\" + `\n  \" It delves into the innermost workings of substances, seeking to understand the nature of atoms, molecules, and their interactions.\" + `\n  \" Chemistry plays a fundamental role in shaping the world we inhabit, influencing everything from materials science to life processes.\" + `\n  \" It empowers us to understand the delicate balance of our planet and address critical challenges such as pollution and climate change.\" + `\n  \" By unraveling the secrets of matter, chemistry becomes an invaluable tool in advancing human knowledge and shaping a sustainable future.\"\nSet-ParagraphText $d.Paragraphs.Item(7) $summaryText\n\n# Append a new empty paragraph at the very end of the document.\n$d.Content.InsertParagraphAfter()\n"}
